$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting the existing rows 51-60 down to 52-61
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly data point
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 45215
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 300000000
$ws.Cells.Item(51, 7).Value = "Espárragos"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 450
$ws.Cells.Item(51, 11).Value = 1300
$ws.Cells.Item(51, 12).Value = 1500
$ws.Cells.Item(51, 13).Value = 1389
$ws.Cells.Item(51, 14).Value = "$/kilo"
$ws.Cells.Item(51, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(51, 16).Value = 1389
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
